# #5: property aircraft done
# Fix property_category values that were left as "land" for the
# "building" (建物) and "car" (汽車) sheets.

$wb = $excel.ActiveWorkbook

# 建物 (building) sheet: column I = property_category, rows 2-6 should read "building"
$wsBuilding = $wb.Worksheets.Item("建物")
for ($r = 2; $r -le 6; $r++) {
    $wsBuilding.Cells.Item($r, 9).Value = "building"
}

# 汽車 (car) sheet: column H = property_category, row 2 should read "car"
$wsCar = $wb.Worksheets.Item("汽車")
$wsCar.Cells.Item(2, 8).Value = "car"
